$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '70.925.46'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +5.79%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.782.48'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +22.53%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '613.27'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +7.26%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '180.07'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.68%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.778.05'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +22.50%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.544'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +6.25%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.168'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +10.78%  '
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.503'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.95%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '41.03'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +14.48%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000259'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +7.78%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.416.26'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +22.55%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.791.52'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +22.83%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '71.033.02'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +5.98%  '
$ws.Range('E18').Value = '  +1.56%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.63'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +9.06%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '525.46'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +8.03%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.80'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.20%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.44'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +23.16%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.750'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +9.84%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '88.71'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +6.53%  '
$ws.Range('E25').Value = '  +10.72%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '13.68'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +8.14%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.99'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +8.14%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0000124'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +31.19%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.52'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +10.32%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.92'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +13.48%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.06'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.76%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '32.35'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +15.69%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.116'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.33%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  +12.42%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.20'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +11.70%  '
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.343'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +10.37%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.23'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +11.16%  '
$ws.Range('E40').Value = '  +7.62%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '51.41'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +4.94%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '431.57'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +17.30%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.159.83'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +13.22%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.90'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +8.39%  '
$ws.Range('E45').Value = '  -6.94%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.81'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.08%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0369'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +7.69%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '27.89'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +9.41%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '141.43'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +5.33%  '
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.49'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +7.82%  '
